$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Send PATCH/POST request for lookup value ..." scenario texts
# (remove "JSON object" suffix, row 26 keeps a trailing space, row 31 does not)
$ws.Range("H15").Value = "Send PATCH request for lookup value"
$ws.Range("H21").Value = "Send PATCH request for lookup value"
$ws.Range("H26").Value = "Send POST request for lookup value "
$ws.Range("H31").Value = "Send POST request for lookup value"

# Row 26 no longer needs the taller wrapped height now that the text is shorter
$ws.Rows.Item(26).AutoFit()

# Update the active cell selection to reflect where editing left off
$null = $ws.Range("I26").Select()
